# Apply the edit described by the diff:
#  1. Slide 22, shape "สี่เหลี่ยมผืนผ้า 4" (the git-commands rectangle): the
#     5th paragraph's text " git push" becomes "git push" (leading space
#     removed, "git" run now leads the paragraph).
#  2. Slides 23 and 24 (the "git checkout -b presentation / create
#     powerpoint" slide and the "edit powerpoint and save as
#     presentation.ppt" slide) are removed entirely from the deck.

$p = $ppt.ActivePresentation

# --- 1. Fix the paragraph text on slide 22 -------------------------------
$s22  = $p.Slides.Item(22)
$shp  = $s22.Shapes.Item(3)
# Captured from the untouched shape (round-tripping the live getter value
# back through the setter loses a bit of precision because Height is a
# Single), so the literal below is what reproduces the original EMU
# extent (2677656) after the autofit box relayouts on edit.
$origHeight = 210.8391

$tr    = $shp.TextFrame.TextRange
$para5 = $tr.Paragraphs(5)

# Before: " git push"  (runs: " ", "git", " push")
# After:  "git push"   (runs: "git", " ", "push")
$para5.Characters(1, 1).Text = ""          # drop the leading space
$para5.Characters(4, 1).Text = " "         # re-touch the inner space so it
                                            # stays a run of its own, giving
                                            # "git" / " " / "push"

# The shape auto-fits its text box; restore its original height so the
# edit doesn't also resize the rectangle.
$shp.Height = $origHeight

# --- 2. Remove the trailing two slides -----------------------------------
$p.Slides.Item(24).Delete()
$p.Slides.Item(23).Delete()
